$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 149, shifting existing rows 149:167 down to 150:168.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new data entry.
$ws.Cells.Item(149, 1).Value = 8
$ws.Cells.Item(149, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44491
$ws.Cells.Item(149, 5).Value = 4
$ws.Cells.Item(149, 6).Value = 100112012
$ws.Cells.Item(149, 7).Value = "Espinaca"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 3400
$ws.Cells.Item(149, 11).Value = 400
$ws.Cells.Item(149, 12).Value = 500
$ws.Cells.Item(149, 13).Value = 450
$ws.Cells.Item(149, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(149, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(149, 16).Value = 900
$ws.Cells.Item(149, 17).Value = 0.5
$ws.Cells.Item(149, 18).Value = "Hortaliza"
